$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Remove the four now-redundant empty "Body Text" spacer paragraphs
# (original 1-based paragraph indices: 11, 15, 72, 93). Deleting from
# the highest index down keeps the earlier indices valid.
# ------------------------------------------------------------------
$toDelete = @(93, 72, 15, 11)
foreach ($i in $toDelete) {
    $d.Paragraphs.Item($i).Range.Delete()
}

# ------------------------------------------------------------------
# Turn the remaining empty "Body Text" paragraph right before the
# "Contacts" bookmark end into a horizontal-rule paragraph (the same
# markup markdown's "---" produces): a run containing a VML rect
# flagged as a horizontal rule.
# ------------------------------------------------------------------
$hrPara = $d.Paragraphs.Item(92)
$hrRange = $hrPara.Range
$hrXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office"><w:body><w:p><w:r><w:pict><v:rect style="width:0;height:1.5pt" o:hralign="center" o:hrstd="t" o:hr="&#116;"/></w:pict></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$hrRange.InsertXML($hrXml)

# ------------------------------------------------------------------
# The "Feedback" intro paragraph switches from the "Block Text" style
# to "First Paragraph".
# ------------------------------------------------------------------
$fbPara = $d.Paragraphs.Item(94)
$fbPara.Style = "First Paragraph"
